# Translated buttons, fixed instructions, skipphases
#
# The worksheet contains three parallel "instruction" blocks laid out side by
# side (columns A-C, D-F and G-I). The third block (G-I, the "Extinction"
# phase) had an extra/incorrect row (a stray "take a break" message) in row 3
# which doesn't belong in that sequence. This shifts all of the remaining
# rows in that block up by one, so row 3 now mirrors the "Again, you will
# see..." message used by the other blocks, and row 7 (which previously held
# the last shifted-in values) becomes empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Extinction" phase instruction column (G, H, I), rows 3-7 ---
# Row 3: replace the stray "break" message with the correct continuation text
$ws.Range("G3").Value = "Again, you will see different faces appear on the screen. Please carefully observe them."
$ws.Range("H3").Value = "同样，你会看到屏幕上出现不同的面孔。请仔细观察它们。"

# Row 4: shift up from what used to be row 5's content
$ws.Range("G4").Value = "When you see the face below, please left-click the mouse as fast as you can"
$ws.Range("H4").Value = "当你看到下面的面孔时，以最快的速度按下空格。"
$ws.Range("I4").Value = "Stimuli/Raw_Trig.BMP"

# Row 5: shift up from what used to be row 6's content
$ws.Range("G5").Value = "When you see any of the faces below, do nothing."
$ws.Range("H5").Value = "当你看到下面的任何一张面孔时，什么都不要做。"
$ws.Range("I5").Value = "Stimuli/Raw_3Faces.BMP"

# Row 6: shift up from what used to be row 7's content
$ws.Range("G6").Value = "Once you are ready, press CONTINUE"
$ws.Range("H6").Value = "一旦你准备好了，按空格"
$ws.Range("I6").Value = "Stimuli/Raw_black.BMP"

# Row 7: now empty (content moved up into row 6)
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("H7").Select()
